$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("TextBox 3")

# Reposition / resize the textbox
# (nudged by < 1 EMU in point-space so float32 round-trip lands on the
# exact target EMU instead of truncating a notch low)
$shp.Left = 209.64712598425197
$shp.Top = 489.1078346456693
$shp.Width = 750.3529527559056
$shp.Height = 50.892244094488184

$tf = $shp.TextFrame
$tf.WordWrap = $true

$tr = $tf.TextRange
$tr.Text = "Reginald Johnson – 19 February 2019`rPresentation and example code can be found at https://github.com/reggie3/testing-examples"
$tr.ParagraphFormat.Alignment = 3

$linkRange = $tr.Characters(83, 43)
$linkRange.ActionSettings.Item(1).Hyperlink.Address = "https://github.com/reggie3/testing-examples"
